$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.338.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.523.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.49%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'597.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.63%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'173.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +3.38%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.134"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +7.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.06%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.135.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.09%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.64%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'67.289.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.51%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.523.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.00%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'14.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'397.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'73.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.89%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.18%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -3.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'10.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.91%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D29").Value = "'6.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.69%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.66%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'24.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'7.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +3.19%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'163.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.64%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'4.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'27.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.86%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'26.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.88%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.801.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.67%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'42.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'341.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.31%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'33.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'6.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.78%  "
$ws.Range("E51").Style = "Normal"
